$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new timesheet entry in row 5: date (A5) + hours worked as a time value (B5),
# matching the formatting already used by the rows above.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 41555

$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 0.17361111111111113

# Remove the (stray) underline that was applied to the hours column, rows 2-5.
$ws.Range("B2:B5").Font.Underline = 0

# Leave the selection where the user ended up after typing the new row.
$ws.Range("C5").Select()
